$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 7899.5
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 7899.5
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 23698.5
$ws.Range("M46").Value = ""
$ws.Range("N46").Value = -23936.5

$ws.Range("H60").Value = 7899.5
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 7899.5
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 23698.5
$ws.Range("M60").Value = ""
$ws.Range("N60").Value = -24666.5

$ws.Range("H98").Value = 2270.182
$ws.Range("I98").Value = 1902.5714
$ws.Range("J98").Value = 9990
$ws.Range("K98").Value = 1902.5714
$ws.Range("L98").Value = 9990
$ws.Range("M98").Value = -404.5714
$ws.Range("N98").Value = -12986

$ws.Range("H103").Value = 863.46155
$ws.Range("I103").Value = 727
$ws.Range("J103").Value = 1318.3334
$ws.Range("K103").Value = 2181
$ws.Range("L103").Value = 3955.0002
$ws.Range("M103").Value = -1595
$ws.Range("N103").Value = -5127.0002

$ws.Range("H105").Value = 7000
$ws.Range("J105").Value = 7000
$ws.Range("L105").Value = 7000
$ws.Range("N105").Value = -13988

$ws.Range("H122").Value = 2270.182
$ws.Range("I122").Value = 1902.5714
$ws.Range("J122").Value = 9990
$ws.Range("K122").Value = 5707.7142
$ws.Range("L122").Value = 29970
$ws.Range("M122").Value = -3257.7142
$ws.Range("N122").Value = -34870

$ws.Range("H132").Value = 1050.5834
$ws.Range("I132").Value = 900.6667
$ws.Range("K132").Value = 2702.0001
$ws.Range("M132").Value = -172.0001000000002

$ws.Range("H137").Value = 16669253
$ws.Range("J137").Value = 2899.9167
$ws.Range("L137").Value = 8699.750100000001
$ws.Range("N137").Value = -13799.7501

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18702.738
$ws.Range("I32").Value = 19294.018
$ws.Range("K32").Value = 19294.018
$ws.Range("M32").Value = -19007.018

$ws.Range("H74").Value = 4837.952
$ws.Range("I74").Value = 3933
$ws.Range("J74").Value = 4988.778
$ws.Range("K74").Value = 3933
$ws.Range("L74").Value = 4988.778
$ws.Range("M74").Value = -3059
$ws.Range("N74").Value = -6736.778

$ws.Range("H77").Value = 4837.952
$ws.Range("I77").Value = 3933
$ws.Range("J77").Value = 4988.778
$ws.Range("K77").Value = 19665
$ws.Range("L77").Value = 24943.89
$ws.Range("M77").Value = -15297
$ws.Range("N77").Value = -33679.89

$ws.Range("H122").Value = 1951.125
$ws.Range("I122").Value = 1951.125
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5853.375
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3403.375
$ws.Range("N122").Value = ""

$ws.Range("H132").Value = 2662.8542
$ws.Range("I132").Value = 2006.55
$ws.Range("K132").Value = 6019.65
$ws.Range("M132").Value = -3489.65

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 3563.875
$ws.Range("J64").Value = 3645
$ws.Range("L64").Value = 3645
$ws.Range("N64").Value = -4095

$ws.Range("H67").Value = 3563.875
$ws.Range("J67").Value = 3645
$ws.Range("L67").Value = 3645
$ws.Range("N67").Value = -5205

$ws.Range("H86").Value = 367713.9
$ws.Range("I86").Value = 4296
$ws.Range("J86").Value = 803815.4
$ws.Range("K86").Value = 4296
$ws.Range("L86").Value = 803815.4
$ws.Range("M86").Value = -3173
$ws.Range("N86").Value = -806061.4

$ws.Range("H89").Value = 367713.9
$ws.Range("I89").Value = 4296
$ws.Range("J89").Value = 803815.4
$ws.Range("K89").Value = 21480
$ws.Range("L89").Value = 4019077
$ws.Range("M89").Value = -15864
$ws.Range("N89").Value = -4030309

$ws.Range("H107").Value = 3710.3333
$ws.Range("I107").Value = 3714.4546
$ws.Range("K107").Value = 3714.4546
$ws.Range("M107").Value = -1794.4546

$ws.Range("H134").Value = 4246.143
$ws.Range("I134").Value = 2070.3333
$ws.Range("K134").Value = 6210.999899999999
$ws.Range("M134").Value = -3675.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4940.3335
$ws.Range("I58").Value = 3042.3333
$ws.Range("J58").Value = 7217.933
$ws.Range("K58").Value = 3042.3333
$ws.Range("L58").Value = 7217.933
$ws.Range("M58").Value = -2839.3333
$ws.Range("N58").Value = -7623.933

$ws.Range("H86").Value = 10050.3
$ws.Range("J86").Value = 12001
$ws.Range("L86").Value = 12001
$ws.Range("N86").Value = -14247

$ws.Range("H89").Value = 10050.3
$ws.Range("J89").Value = 12001
$ws.Range("L89").Value = 60005
$ws.Range("N89").Value = -71237

$ws.Range("H134").Value = 8387.643
$ws.Range("I134").Value = 7913.615
$ws.Range("K134").Value = 23740.845
$ws.Range("M134").Value = -21205.845

$ws.Range("H136").Value = 4940.3335
$ws.Range("I136").Value = 3042.3333
$ws.Range("J136").Value = 7217.933
$ws.Range("K136").Value = 9126.999899999999
$ws.Range("L136").Value = 21653.799
$ws.Range("M136").Value = -6576.999899999999
$ws.Range("N136").Value = -26753.799

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 687.25
$ws.Range("I63").Value = 249.66667
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 749.00001
$ws.Range("L63").Value = 6000
$ws.Range("M63").Value = -0.000009999999974752427
$ws.Range("N63").Value = -7498

$ws.Range("H64").Value = 10013
$ws.Range("J64").Value = 10014
$ws.Range("L64").Value = 30042
$ws.Range("N64").Value = -30582

$ws.Range("H66").Value = 687.25
$ws.Range("I66").Value = 249.66667
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 2247.00003
$ws.Range("L66").Value = 18000
$ws.Range("M66").Value = 1496.99997
$ws.Range("N66").Value = -25488

$ws.Range("H67").Value = 10013
$ws.Range("J67").Value = 10014
$ws.Range("L67").Value = 30042
$ws.Range("N67").Value = -31914

$ws.Range("H74").Value = 11669.333
$ws.Range("J74").Value = 11997.5
$ws.Range("L74").Value = 35992.5
$ws.Range("N74").Value = -38114.5

$ws.Range("H77").Value = 11669.333
$ws.Range("J77").Value = 11997.5
$ws.Range("L77").Value = 107977.5
$ws.Range("N77").Value = -118585.5

$ws.Range("H101").Value = 7674.75
$ws.Range("I101").Value = 8500
$ws.Range("J101").Value = 6849.5
$ws.Range("K101").Value = 25500
$ws.Range("L101").Value = 20548.5
$ws.Range("M101").Value = -23066
$ws.Range("N101").Value = -25416.5

$ws.Range("H104").Value = 6324.875
$ws.Range("I104").Value = 10649.75
$ws.Range("K104").Value = 31949.25
$ws.Range("M104").Value = -29328.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1884.4286
$ws.Range("I102").Value = 1211
$ws.Range("K102").Value = 1211
$ws.Range("M102").Value = 411

$ws.Range("H122").Value = 4740.032
$ws.Range("I122").Value = 4051.8262
$ws.Range("J122").Value = 6718.625
$ws.Range("K122").Value = 12155.4786
$ws.Range("L122").Value = 20155.875
$ws.Range("M122").Value = -9705.4786
$ws.Range("N122").Value = -25055.875

$ws.Range("H126").Value = 2741.75
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 51786.75
$ws.Range("I74").Value = 49049
$ws.Range("J74").Value = 60000
$ws.Range("K74").Value = 49049
$ws.Range("L74").Value = 60000
$ws.Range("M74").Value = -48051
$ws.Range("N74").Value = -61996

$ws.Range("H77").Value = 51786.75
$ws.Range("I77").Value = 49049
$ws.Range("J77").Value = 60000
$ws.Range("K77").Value = 147147
$ws.Range("L77").Value = 180000
$ws.Range("M77").Value = -142155
$ws.Range("N77").Value = -189984

$ws.Range("H82").Value = 2072.8333
$ws.Range("I82").Value = 1516.6666
$ws.Range("J82").Value = 2350.9167
$ws.Range("K82").Value = 1516.6666
$ws.Range("L82").Value = 2350.9167
$ws.Range("M82").Value = -1155.6666
$ws.Range("N82").Value = -3072.9167

$ws.Range("H85").Value = 2072.8333
$ws.Range("I85").Value = 1516.6666
$ws.Range("J85").Value = 2350.9167
$ws.Range("K85").Value = 1516.6666
$ws.Range("L85").Value = 2350.9167
$ws.Range("M85").Value = -268.6666
$ws.Range("N85").Value = -4846.9167

$ws.Range("H136").Value = 4350.0415
$ws.Range("I136").Value = 1915.1818
$ws.Range("K136").Value = 5745.5454
$ws.Range("M136").Value = -3195.5454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1508.05
$ws.Range("I107").Value = 1449.4706
$ws.Range("J107").Value = 1840
$ws.Range("K107").Value = 4348.4118
$ws.Range("L107").Value = 5520
$ws.Range("M107").Value = -2428.4118
$ws.Range("N107").Value = -9360

$ws.Range("H126").Value = 2215.5334
$ws.Range("I126").Value = 2215.5334
$ws.Range("K126").Value = 6646.600199999999
$ws.Range("M126").Value = -4176.600199999999
